# Estadisticos Segundo Parcial 26 Mayo
# Updates the 2nd-partial statistics (sheet "Estadisticos 2P"), recomputed
# final statistics (sheet "Estadisticos Final") for group 2ALCV, and
# refreshes the "Rescatables" (make-up exam candidates) listing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Estadisticos 2P" - row 3 (Ingles II / 2ALCV) totals updated and a
#    Promedio (H3) value added.
# ---------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D3").Value = 0
$ws2P.Range("E3").Value = 6
$ws2P.Range("F3").Value = 22
$ws2P.Range("G3").Value = 78.57
$ws2P.Range("H3").Value = 5.7

# ---------------------------------------------------------------------
# 2) "Estadisticos Final" - row 3 (Ingles II / 2ALCV) totals updated to
#    match the new second-partial figures.
# ---------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("E3").Value = 6
$wsFinal.Range("F3").Value = 22
$wsFinal.Range("G3").Value = 78.57
$wsFinal.Range("H3").Value = 6.8

# ---------------------------------------------------------------------
# 3) "Rescatables" - the make-up-exam candidate list is rebuilt: it now
#    has 21 students (rows 2-22) instead of 29 (rows 2-30), so the extra
#    rows are removed and the remaining rows are rewritten in full.
# ---------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

# Remove the rows that no longer exist (old rows 23-30).
$wsResc.Rows("23:30").Delete()

$rescatablesData = @(
    @(24330051920304, 'ARMAS', 'SALINAS', 'JOSE GUSTAVO', 'Ingles II', '2AEV', 4),
    @(24330051920305, 'MORALES', 'CUAHUA', 'ANDRES', 'Ingles II', '2AEV', 4),
    @(24330051920113, 'RAMOS', 'DE LA CRUZ', 'DEREK', 'Ingles II', '2AEV', 4),
    @(24330051920330, 'VASQUEZ', 'PEREZ', 'DANIELA LILI', 'Ingles II', '2ALCV', 4),
    @(24330051920246, 'ZUNO', 'FLORES', 'ALIN MARIEL', 'Ingles II', '2ALCV', 4),
    @(24330051920393, 'MUÑOZ', 'REYES', 'ERWIN ISRAEL', 'Ingles II', '2APV', 4),
    @(23330051920211, 'VAZQUEZ', 'CARRILLO', 'DIEGO ARMANDO', 'Ingles IV', '4AEV', 4),
    @(24330051920093, 'ARIAS', 'SARMIENTO', 'URIEL ARTURO', 'Ingles II', '2AEV', 3),
    @(24330051920098, 'CHICO', 'BALDERAS', 'YARETH', 'Ingles II', '2AEV', 3),
    @(24330051920144, 'MUÑOZ', 'CORONA', 'JOSE ABEL', 'Ingles II', '2AEV', 3),
    @(24330051920143, 'ROSAS', 'MEZA', 'CARLOS ANTONIO', 'Ingles II', '2AEV', 3),
    @(22330051920389, 'FLORES', 'LAGUNA', 'JOSE ANTONIO', 'Ingles IV', '4AEV', 3),
    @(23330051920332, 'RODRIGUEZ', 'SUAREZ', 'SERGIO JOSUE', 'Ingles IV', '4AEV', 3),
    @(24330051920392, 'CERON', 'GONZALEZ', 'LEVI SANTIAGO', 'Ingles II', '2AEV', 2),
    @(23330051920224, 'DORANTES', 'PORRAS', 'ROBERTO', 'Ingles II', '2AEV', 2),
    @(24330051920220, 'GARCIA', 'CHAPARRO', 'MAYKA XIMENA', 'Ingles II', '2ALCV', 2),
    @(24330051920226, 'LEYVA', 'HERNANDEZ', 'EUNICE GUADALUPE', 'Ingles II', '2ALCV', 2),
    @(23330051920329, 'JIMENEZ', 'CIRUELO', 'SABDY', 'Ingles IV', '4AEV', 2),
    @(23330051920324, 'JUAREZ', 'LIBRADO', 'ARMANDO GABRIEL', 'Ingles IV', '4AEV', 2),
    @(23330051920203, 'PEREZ', 'DE JESUS', 'LUIS FABIAN', 'Ingles IV', '4AEV', 2),
    @(24330051920396, 'MARTINEZ', 'GONZALEZ', 'SANTIAGO', 'Ingles II', '2ALCV', 1)
)

$rowCount = $rescatablesData.Count
$colCount = 7
$arr = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i, $j] = $rescatablesData[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$targetRange = $wsResc.Range($wsResc.Cells.Item($startRow, 1), $wsResc.Cells.Item($endRow, 7))
$targetRange.Value = $arr
